$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 ---------------------------------------------------------------
$ws.Range("B13").Value = 164
$ws.Range("B14").Value = 167

$ws.Range("C13").Value = " You\'ve graduated!\nCongratulations!"
$ws.Range("C14").Value = " We hope for your continued\npatronage and generosity! ♪"

$ws.Range("D13").Value = " Вы выпускники! Поздравляю!"
$ws.Range("D14").Value = " Мы надеемся на ваше дальнейшее\nпокровительство и щедрость! ♪"

$ws.Range("E13").Value = " Âú âúðôòëîéëé! Ðïèäñàâìÿý!"
$ws.Range("E14").Value = " Íú îàäååíòÿ îà âàšå äàìûîåêšåå\nðïëñïâéóåìûòóâï é þåäñïòóû! ♪"

# Row heights (row 13 was manually resized in the source; row 14 matches the
# auto-fit height already used elsewhere in the sheet for similarly sized
# two-line cells).
$ws.Rows.Item(13).RowHeight = 22.8
$ws.Rows.Item(14).RowHeight = 31.8

# --- Selection / view state ------------------------------------------------
$ws.Range("D16").Select()
